$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 2015.1111
$ws.Range("J40").Value2 = 2019.4286
$ws.Range("L40").Value2 = 2019.4286
$ws.Range("N40").Value2 = -2369.4286
$ws.Range("H138").Value2 = 288435.47
$ws.Range("I138").Value2 = 440143.1
$ws.Range("K138").Value2 = 1320429.3
$ws.Range("M138").Value2 = -1315289.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3320.29
$ws.Range("I32").Value2 = 3308.3738
$ws.Range("K32").Value2 = 3308.3738
$ws.Range("M32").Value2 = -3021.3738
$ws.Range("H35").Value2 = 5999.8
$ws.Range("I35").Value2 = 3333
$ws.Range("J35").Value2 = 10000
$ws.Range("K35").Value2 = 3333
$ws.Range("L35").Value2 = 10000
$ws.Range("M35").Value2 = -2927
$ws.Range("N35").Value2 = -10812
$ws.Range("H74").Value2 = 3545.818
$ws.Range("J74").Value2 = 6197.375
$ws.Range("L74").Value2 = 6197.375
$ws.Range("N74").Value2 = -7945.375
$ws.Range("H77").Value2 = 3545.818
$ws.Range("J77").Value2 = 6197.375
$ws.Range("L77").Value2 = 30986.875
$ws.Range("N77").Value2 = -39722.875
$ws.Range("H109").Value2 = 0
$ws.Range("J109").Value2 = 0
$ws.Range("L109").Value2 = 0
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value2 = 3372.4
$ws.Range("I110").Value2 = 2841.3635
$ws.Range("K110").Value2 = 2841.3635
$ws.Range("M110").Value2 = -796.3634999999999
$ws.Range("H122").Value2 = 1261143.8
$ws.Range("I122").Value2 = 3764.7
$ws.Range("J122").Value2 = 2937649.2
$ws.Range("K122").Value2 = 11294.1
$ws.Range("L122").Value2 = 8812947.600000001
$ws.Range("M122").Value2 = -8844.099999999999
$ws.Range("N122").Value2 = -8817847.600000001
$ws.Range("H135").Value2 = 125843.25
$ws.Range("J135").Value2 = 125843.25
$ws.Range("L135").Value2 = 125843.25
$ws.Range("N135").Value2 = -135983.25
$ws.Range("H139").Value2 = 39999.5
$ws.Range("J139").Value2 = 39999.5
$ws.Range("L139").Value2 = 39999.5
$ws.Range("N139").Value2 = -50279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value2 = 21999.75
$ws.Range("I75").Value2 = 15666.333
$ws.Range("K75").Value2 = 15666.333
$ws.Range("M75").Value2 = -14730.333
$ws.Range("H78").Value2 = 21999.75
$ws.Range("I78").Value2 = 15666.333
$ws.Range("K78").Value2 = 46998.999
$ws.Range("M78").Value2 = -42318.999
$ws.Range("H81").Value2 = 18942.5
$ws.Range("J81").Value2 = 18942.5
$ws.Range("L81").Value2 = 18942.5
$ws.Range("N81").Value2 = -21064.5
$ws.Range("H84").Value2 = 18942.5
$ws.Range("J84").Value2 = 18942.5
$ws.Range("L84").Value2 = 56827.5
$ws.Range("N84").Value2 = -67435.5
$ws.Range("H94").Value2 = 7196.317
$ws.Range("I94").Value2 = 7860.886
$ws.Range("J94").Value2 = 5368.75
$ws.Range("K94").Value2 = 7860.886
$ws.Range("L94").Value2 = 5368.75
$ws.Range("M94").Value2 = -7409.886
$ws.Range("N94").Value2 = -6270.75
$ws.Range("H100").Value2 = 36399.6
$ws.Range("J100").Value2 = 36399.6
$ws.Range("L100").Value2 = 36399.6
$ws.Range("N100").Value2 = -38563.6
$ws.Range("H135").Value2 = 84497
$ws.Range("J135").Value2 = 89329.664
$ws.Range("L135").Value2 = 89329.664
$ws.Range("N135").Value2 = -99469.664
$ws.Range("H139").Value2 = 60000
$ws.Range("I139").Value2 = 40000
$ws.Range("J139").Value2 = 80000
$ws.Range("K139").Value2 = 40000
$ws.Range("L139").Value2 = 80000
$ws.Range("M139").Value2 = -34860
$ws.Range("N139").Value2 = -90280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1199.0834
$ws.Range("J22").Value2 = 1774.75
$ws.Range("L22").Value2 = 1774.75
$ws.Range("N22").Value2 = -2474.75
$ws.Range("H31").Value2 = 3433.6135
$ws.Range("I31").Value2 = 2839.919
$ws.Range("K31").Value2 = 2839.919
$ws.Range("M31").Value2 = -2544.919
$ws.Range("H34").Value2 = 3433.6135
$ws.Range("I34").Value2 = 2839.919
$ws.Range("K34").Value2 = 2839.919
$ws.Range("M34").Value2 = -2637.919
$ws.Range("H58").Value2 = 1031.2307
$ws.Range("I58").Value2 = 700.9167
$ws.Range("K58").Value2 = 700.9167
$ws.Range("M58").Value2 = -497.9167
$ws.Range("H59").Value2 = 64563.5
$ws.Range("J59").Value2 = 64563.5
$ws.Range("L59").Value2 = 64563.5
$ws.Range("N59").Value2 = -66853.5
$ws.Range("H60").Value2 = 11999
$ws.Range("J60").Value2 = 14998.5
$ws.Range("L60").Value2 = 14998.5
$ws.Range("N60").Value2 = -16020.5
$ws.Range("H99").Value2 = 299884.8
$ws.Range("I99").Value2 = 719113.5600000001
$ws.Range("K99").Value2 = 719113.5600000001
$ws.Range("M99").Value2 = -717615.5600000001
$ws.Range("H122").Value2 = 4783.892
$ws.Range("I122").Value2 = 6168
$ws.Range("K122").Value2 = 18504
$ws.Range("M122").Value2 = -16054
$ws.Range("H126").Value2 = 299884.8
$ws.Range("I126").Value2 = 719113.5600000001
$ws.Range("K126").Value2 = 2157340.68
$ws.Range("M126").Value2 = -2154870.68
$ws.Range("H136").Value2 = 1031.2307
$ws.Range("I136").Value2 = 700.9167
$ws.Range("K136").Value2 = 2102.7501
$ws.Range("M136").Value2 = 447.2498999999998
$ws.Range("H141").Value2 = 203495.67
$ws.Range("I141").Value2 = 0
$ws.Range("K141").Value2 = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 911998.8
$ws.Range("I5").Value2 = 1847.25
$ws.Range("K5").Value2 = 5541.75
$ws.Range("M5").Value2 = -5429.75
$ws.Range("H19").Value2 = 20
$ws.Range("J19").Value2 = 0
$ws.Range("L19").Value2 = 0
$ws.Range("N19").ClearContents()
$ws.Range("H60").Value2 = 921.3333
$ws.Range("I60").Value2 = 111.875
$ws.Range("K60").Value2 = 335.625
$ws.Range("M60").Value2 = -84.625
$ws.Range("H107").Value2 = 848.1923
$ws.Range("I107").Value2 = 302.25
$ws.Range("J107").Value2 = 947.4545000000001
$ws.Range("K107").Value2 = 906.75
$ws.Range("L107").Value2 = 2842.3635
$ws.Range("M107").Value2 = 1013.25
$ws.Range("N107").Value2 = -6682.3635
$ws.Range("H114").Value2 = 1711.6428
$ws.Range("J114").Value2 = 2239.375
$ws.Range("L114").Value2 = 6718.125
$ws.Range("N114").Value2 = -13226.125
$ws.Range("H115").Value2 = 2498.1428
$ws.Range("J115").Value2 = 799
$ws.Range("L115").Value2 = 2397
$ws.Range("N115").Value2 = -4747
$ws.Range("H135").Value2 = 911998.8
$ws.Range("I135").Value2 = 1847.25
$ws.Range("K135").Value2 = 16625.25
$ws.Range("M135").Value2 = -14090.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 8655.462
$ws.Range("I122").Value2 = 6191.35
$ws.Range("J122").Value2 = 16869.166
$ws.Range("K122").Value2 = 18574.05
$ws.Range("L122").Value2 = 50607.49800000001
$ws.Range("M122").Value2 = -16124.05
$ws.Range("N122").Value2 = -55507.49800000001
$ws.Range("H123").Value2 = 38333.332
$ws.Range("J123").Value2 = 38333.332
$ws.Range("L123").Value2 = 38333.332
$ws.Range("N123").Value2 = -43233.332
$ws.Range("H132").Value2 = 4108.885
$ws.Range("I132").Value2 = 3820.8718
$ws.Range("J132").Value2 = 4972.923
$ws.Range("K132").Value2 = 11462.6154
$ws.Range("L132").Value2 = 14918.769
$ws.Range("M132").Value2 = -8932.615399999999
$ws.Range("N132").Value2 = -19978.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 15459.383
$ws.Range("I7").Value2 = 22036.26
$ws.Range("K7").Value2 = 22036.26
$ws.Range("M7").Value2 = -21924.26
$ws.Range("H16").Value2 = 2553.2273
$ws.Range("I16").Value2 = 2889.1052
$ws.Range("J16").Value2 = 426
$ws.Range("K16").Value2 = 2889.1052
$ws.Range("L16").Value2 = 426
$ws.Range("M16").Value2 = -2719.1052
$ws.Range("N16").Value2 = -766
$ws.Range("H22").Value2 = 14066.087
$ws.Range("I22").Value2 = 20023.77
$ws.Range("J22").Value2 = 6321.1
$ws.Range("K22").Value2 = 20023.77
$ws.Range("L22").Value2 = 6321.1
$ws.Range("M22").Value2 = -19728.77
$ws.Range("N22").Value2 = -6911.1
$ws.Range("H27").Value2 = 14066.087
$ws.Range("I27").Value2 = 20023.77
$ws.Range("J27").Value2 = 6321.1
$ws.Range("K27").Value2 = 20023.77
$ws.Range("L27").Value2 = 6321.1
$ws.Range("M27").Value2 = -19916.77
$ws.Range("N27").Value2 = -6535.1
$ws.Range("H40").Value2 = 25988.85
$ws.Range("I40").Value2 = 47254.445
$ws.Range("J40").Value2 = 8589.727999999999
$ws.Range("K40").Value2 = 47254.445
$ws.Range("L40").Value2 = 8589.727999999999
$ws.Range("M40").Value2 = -47118.445
$ws.Range("N40").Value2 = -8861.727999999999
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 11085.429
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 11085.429
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -11461.429
$ws.Range("H100").Value2 = 2416
$ws.Range("I100").Value2 = 2416
$ws.Range("J100").Value2 = 0
$ws.Range("K100").Value2 = 2416
$ws.Range("L100").Value2 = 0
$ws.Range("M100").Value2 = -1875
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value2 = 9619.134
$ws.Range("I122").Value2 = 9619.134
$ws.Range("K122").Value2 = 28857.402
$ws.Range("M122").Value2 = -26407.402
$ws.Range("H126").Value2 = 15459.383
$ws.Range("I126").Value2 = 22036.26
$ws.Range("K126").Value2 = 66108.78
$ws.Range("M126").Value2 = -63638.78
$ws.Range("H132").Value2 = 295494.62
$ws.Range("I132").Value2 = 467765.4
$ws.Range("K132").Value2 = 1403296.2
$ws.Range("M132").Value2 = -1400766.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 10353.851
$ws.Range("I132").Value2 = 11038.771
$ws.Range("J132").Value2 = 6449.8
$ws.Range("K132").Value2 = 33116.313
$ws.Range("L132").Value2 = 19349.4
$ws.Range("M132").Value2 = -30586.313
$ws.Range("N132").Value2 = -24409.4
